$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C2:C4 with the Tsai Model formula (10^6)
$ws.Range("C2").Formula = "=10^6"
$ws.Range("C3").Formula = "=10^6"
$ws.Range("C4").Formula = "=10^6"

# Update the selected/active cell to C3
$ws.Range("C3").Select()
